$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 19659
$ws.Range("E2").Value = 782
$ws.Range("F2").Value = 782
$ws.Range("G2").Value = 1286
$ws.Range("H2").Value = 979
$ws.Range("I2").Value = 976
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 18830
$ws.Range("L2").Value = 5012
$ws.Range("M2").Value = 13817
$ws.Range("N2").Value = 13767
$ws.Range("O2").Value = 50
$ws.Range("P2").Value = 489
$ws.Range("Q2").Value = 624
$ws.Range("R2").Value = -479
$ws.Range("S2").Value = -170
$ws.Range("T2").Value = 185
$ws.Range("U2").Value = 439
$ws.Range("V2").Value = 155
$ws.Range("W2").Value = 3.98
$ws.Range("X2").Value = 4.98
$ws.Range("Y2").Value = 6.99
$ws.Range("Z2").Value = 5.15
$ws.Range("AA2").Value = 36.28
$ws.Range("AB2").Value = 2343.27
$ws.Range("AC2").Value = 999
$ws.Range("AD2").Value = 18.97
$ws.Range("AE2").Value = 15771
$ws.Range("AF2").Value = 1.2
$ws.Range("AG2").Value = 60
$ws.Range("AH2").Value = 0.32
$ws.Range("AI2").Value = 5.37
$ws.Range("AJ2").Value = 97704482

# Row 3
$ws.Range("D3").Value = 21128
$ws.Range("E3").Value = 879
$ws.Range("F3").Value = 879
$ws.Range("G3").Value = 1206
$ws.Range("H3").Value = 915
$ws.Range("I3").Value = 908
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 19541
$ws.Range("L3").Value = 4862
$ws.Range("M3").Value = 14679
$ws.Range("N3").Value = 14622
$ws.Range("O3").Value = 57
$ws.Range("P3").Value = 489
$ws.Range("Q3").Value = 540
$ws.Range("R3").Value = -458
$ws.Range("S3").Value = -93
$ws.Range("T3").Value = 290
$ws.Range("U3").Value = 250
$ws.Range("V3").Value = 115
$ws.Range("W3").Value = 4.16
$ws.Range("X3").Value = 4.33
$ws.Range("Y3").Value = 6.4
$ws.Range("Z3").Value = 4.77
$ws.Range("AA3").Value = 33.13
$ws.Range("AB3").Value = 2513.04
$ws.Range("AC3").Value = 929
$ws.Range("AD3").Value = 26.68
$ws.Range("AE3").Value = 16750
$ws.Range("AF3").Value = 1.48
$ws.Range("AG3").Value = 60
$ws.Range("AH3").Value = 0.24
$ws.Range("AI3").Value = 5.77
$ws.Range("AJ3").Value = 97704482

# Row 4
$ws.Range("D4").Value = 25217
$ws.Range("E4").Value = 1052
$ws.Range("F4").Value = 1052
$ws.Range("G4").Value = 1387
$ws.Range("H4").Value = 1053
$ws.Range("I4").Value = 949
$ws.Range("J4").Value = 104
$ws.Range("K4").Value = 21817
$ws.Range("L4").Value = 5747
$ws.Range("M4").Value = 16071
$ws.Range("N4").Value = 15096
$ws.Range("O4").Value = 975
$ws.Range("P4").Value = 489
$ws.Range("Q4").Value = 1733
$ws.Range("R4").Value = -964
$ws.Range("S4").Value = -666
$ws.Range("T4").Value = 178
$ws.Range("U4").Value = 1556
$ws.Range("V4").Value = 327
$ws.Range("W4").Value = 4.17
$ws.Range("X4").Value = 4.17
$ws.Range("Y4").Value = 6.39
$ws.Range("Z4").Value = 5.09
$ws.Range("AA4").Value = 35.76
$ws.Range("AB4").Value = 2693.86
$ws.Range("AC4").Value = 971
$ws.Range("AD4").Value = 15.86
$ws.Range("AE4").Value = 17293
$ws.Range("AF4").Value = 0.89
$ws.Range("AG4").Value = 60
$ws.Range("AH4").Value = 0.39
$ws.Range("AI4").Value = 5.52
$ws.Range("AJ4").Value = 97704482

# Row 5
$ws.Range("D5").Value = 25340
$ws.Range("E5").Value = 871
$ws.Range("F5").Value = 871
$ws.Range("G5").Value = 1607
$ws.Range("H5").Value = 1215
$ws.Range("I5").Value = 1134
$ws.Range("J5").Value = 81
$ws.Range("K5").Value = 27969
$ws.Range("L5").Value = 7535
$ws.Range("M5").Value = 20434
$ws.Range("N5").Value = 16739
$ws.Range("O5").Value = 3694
$ws.Range("P5").Value = 489
$ws.Range("Q5").Value = 1064
$ws.Range("R5").Value = -829
$ws.Range("S5").Value = -54
$ws.Range("T5").Value = 240
$ws.Range("U5").Value = 825
$ws.Range("V5").Value = 373
$ws.Range("W5").Value = 3.44
$ws.Range("X5").Value = 4.79
$ws.Range("Y5").Value = 7.12
$ws.Range("Z5").Value = 4.88
$ws.Range("AA5").Value = 36.88
$ws.Range("AB5").Value = 3033.2
$ws.Range("AC5").Value = 1161
$ws.Range("AD5").Value = 12.88
$ws.Range("AE5").Value = 19173
$ws.Range("AF5").Value = 0.78
$ws.Range("AG5").Value = 80
$ws.Range("AH5").Value = 0.54
$ws.Range("AI5").Value = 6.16
$ws.Range("AJ5").Value = 97704482

# Row 6
$ws.Range("D6").Value = 32517
$ws.Range("E6").Value = 1372
$ws.Range("F6").Value = 1372
$ws.Range("G6").Value = 1823
$ws.Range("H6").Value = 1272
$ws.Range("I6").Value = 1031
$ws.Range("K6").Value = 28374
$ws.Range("L6").Value = 7263
$ws.Range("M6").Value = 21111
$ws.Range("N6").Value = 17082
$ws.Range("P6").Value = 489
$ws.Range("Q6").Value = 954
$ws.Range("R6").Value = -984
$ws.Range("S6").Value = 86
$ws.Range("T6").Value = 716
$ws.Range("U6").Value = 239
$ws.Range("V6").Value = 348
$ws.Range("W6").Value = 4.22
$ws.Range("X6").Value = 3.91
$ws.Range("Y6").Value = 6.1
$ws.Range("Z6").Value = 4.52
$ws.Range("AA6").Value = 34.4
$ws.Range("AB6").Value = 3231.48
$ws.Range("AC6").Value = 1056
$ws.Range("AD6").Value = 13.55
$ws.Range("AE6").Value = 19565
$ws.Range("AF6").Value = 0.73
$ws.Range("AG6").Value = 210
$ws.Range("AH6").Value = 1.47
$ws.Range("AI6").Value = 17.78
$ws.Range("AJ6").Value = 97704482

# Row 7
$ws.Range("D7").Value = 31186
$ws.Range("E7").Value = 1069
$ws.Range("G7").Value = 1645
$ws.Range("H7").Value = 1214
$ws.Range("I7").Value = 1100
$ws.Range("K7").Value = 29273
$ws.Range("L7").Value = 7630
$ws.Range("M7").Value = 21644
$ws.Range("N7").Value = 17926
$ws.Range("P7").Value = 490
$ws.Range("Q7").Value = 1043
$ws.Range("R7").Value = -935
$ws.Range("S7").Value = -165
$ws.Range("T7").Value = 443
$ws.Range("U7").Value = 20
$ws.Range("W7").Value = 3.43
$ws.Range("X7").Value = 3.89
$ws.Range("Y7").Value = 6.28
$ws.Range("Z7").Value = 4.21
$ws.Range("AA7").Value = 35.25
$ws.Range("AC7").Value = 1125
$ws.Range("AD7").Value = 9.15
$ws.Range("AE7").Value = 20532
$ws.Range("AF7").Value = 0.5
$ws.Range("AG7").Value = 210
$ws.Range("AH7").Value = 2.04
$ws.Range("AI7").Value = 18.66

# Row 8
$ws.Range("D8").Value = 32237
$ws.Range("E8").Value = 1143
$ws.Range("G8").Value = 1729
$ws.Range("H8").Value = 1282
$ws.Range("I8").Value = 1178
$ws.Range("K8").Value = 30480
$ws.Range("L8").Value = 7788
$ws.Range("M8").Value = 22693
$ws.Range("N8").Value = 18862
$ws.Range("P8").Value = 490
$ws.Range("Q8").Value = 1665
$ws.Range("R8").Value = -1858
$ws.Range("S8").Value = 30
$ws.Range("T8").Value = 443
$ws.Range("U8").Value = 460
$ws.Range("W8").Value = 3.55
$ws.Range("X8").Value = 3.98
$ws.Range("Y8").Value = 6.4
$ws.Range("Z8").Value = 4.29
$ws.Range("AA8").Value = 34.32
$ws.Range("AC8").Value = 1205
$ws.Range("AD8").Value = 8.550000000000001
$ws.Range("AE8").Value = 21604
$ws.Range("AF8").Value = 0.48
$ws.Range("AG8").Value = 210
$ws.Range("AH8").Value = 2.04
$ws.Range("AI8").Value = 17.42

# Row 9
$ws.Range("D9").Value = 33493
$ws.Range("E9").Value = 1219
$ws.Range("G9").Value = 1837
$ws.Range("H9").Value = 1364
$ws.Range("I9").Value = 1252
$ws.Range("K9").Value = 31804
$ws.Range("L9").Value = 7983
$ws.Range("M9").Value = 23821
$ws.Range("N9").Value = 19864
$ws.Range("P9").Value = 490
$ws.Range("Q9").Value = 1180
$ws.Range("R9").Value = -1840
$ws.Range("S9").Value = 33
$ws.Range("T9").Value = 450
$ws.Range("U9").Value = 390
$ws.Range("W9").Value = 3.64
$ws.Range("X9").Value = 4.07
$ws.Range("Y9").Value = 6.47
$ws.Range("Z9").Value = 4.38
$ws.Range("AA9").Value = 33.51
$ws.Range("AC9").Value = 1282
$ws.Range("AD9").Value = 8.039999999999999
$ws.Range("AE9").Value = 22752
$ws.Range("AF9").Value = 0.45
$ws.Range("AG9").Value = 210
$ws.Range("AH9").Value = 2.04
$ws.Range("AI9").Value = 16.38
